$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 225.81818
$ws.Range("I53").Value = 194.85715
$ws.Range("J53").Value = 280
$ws.Range("K53").Value = 194.85715
$ws.Range("L53").Value = 280
$ws.Range("M53").Value = 442.14285
$ws.Range("N53").Value = -1554
$ws.Range("H54").Value = 15534.5
$ws.Range("I54").Value = 10000
$ws.Range("J54").Value = 21069
$ws.Range("K54").Value = 10000
$ws.Range("L54").Value = 21069
$ws.Range("M54").Value = -9514
$ws.Range("N54").Value = -22041
$ws.Range("H112").Value = 2096.8667
$ws.Range("I112").Value = 1313.909
$ws.Range("K112").Value = 3941.727
$ws.Range("M112").Value = -2833.727
$ws.Range("H133").Value = 71994.5
$ws.Range("J133").Value = 71994.5
$ws.Range("L133").Value = 71994.5
$ws.Range("N133").Value = -82114.5
$ws.Range("H137").Value = 11732.533
$ws.Range("I137").Value = 1931.6666
$ws.Range("K137").Value = 5794.9998
$ws.Range("M137").Value = -3244.9998
$ws.Range("H138").Value = 2034.9166
$ws.Range("I138").Value = 1339.0294
$ws.Range("J138").Value = 2508.12
$ws.Range("K138").Value = 4017.0882
$ws.Range("L138").Value = 7524.36
$ws.Range("M138").Value = 1122.9118
$ws.Range("N138").Value = -17804.36

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 753.5
$ws.Range("J2").Value = 799
$ws.Range("L2").Value = 799
$ws.Range("N2").Value = -1025
$ws.Range("H32").Value = 112173.14
$ws.Range("I32").Value = 119115.58
$ws.Range("J32").Value = 12664.833
$ws.Range("K32").Value = 119115.58
$ws.Range("L32").Value = 12664.833
$ws.Range("M32").Value = -118828.58
$ws.Range("N32").Value = -13238.833
$ws.Range("H63").Value = 2249.5454
$ws.Range("I63").Value = 2249.611
$ws.Range("J63").Value = 2249.25
$ws.Range("K63").Value = 2249.611
$ws.Range("L63").Value = 2249.25
$ws.Range("M63").Value = -1563.611
$ws.Range("N63").Value = -3621.25
$ws.Range("H66").Value = 2249.5454
$ws.Range("I66").Value = 2249.611
$ws.Range("J66").Value = 2249.25
$ws.Range("K66").Value = 11248.055
$ws.Range("L66").Value = 11246.25
$ws.Range("M66").Value = -7816.055
$ws.Range("N66").Value = -18110.25
$ws.Range("H74").Value = 5118.353
$ws.Range("I74").Value = 916.5814
$ws.Range("K74").Value = 916.5814
$ws.Range("M74").Value = -42.58140000000003
$ws.Range("H77").Value = 5118.353
$ws.Range("I77").Value = 916.5814
$ws.Range("K77").Value = 4582.907
$ws.Range("M77").Value = -214.9070000000002
$ws.Range("H80").Value = 29475.572
$ws.Range("I80").Value = 20750
$ws.Range("J80").Value = 32965.8
$ws.Range("K80").Value = 20750
$ws.Range("L80").Value = 32965.8
$ws.Range("M80").Value = -19752
$ws.Range("N80").Value = -34961.8
$ws.Range("H83").Value = 29475.572
$ws.Range("I83").Value = 20750
$ws.Range("J83").Value = 32965.8
$ws.Range("K83").Value = 62250
$ws.Range("L83").Value = 98897.40000000001
$ws.Range("M83").Value = -57258
$ws.Range("N83").Value = -108881.4
$ws.Range("H88").Value = 3292.5715
$ws.Range("J88").Value = 3429.6
$ws.Range("L88").Value = 3429.6
$ws.Range("N88").Value = -4241.6
$ws.Range("H91").Value = 3292.5715
$ws.Range("J91").Value = 3429.6
$ws.Range("L91").Value = 3429.6
$ws.Range("N91").Value = -6237.6
$ws.Range("H97").Value = 32259092
$ws.Range("I97").Value = 753.9655
$ws.Range("K97").Value = 753.9655
$ws.Range("M97").Value = -257.9655
$ws.Range("H116").Value = 753.5
$ws.Range("J116").Value = 799
$ws.Range("L116").Value = 799
$ws.Range("N116").Value = -5387

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 753.5
$ws.Range("J3").Value = 799
$ws.Range("L3").Value = 799
$ws.Range("N3").Value = -1027
$ws.Range("H132").Value = 99990
$ws.Range("J132").Value = 99990
$ws.Range("L132").Value = 99990
$ws.Range("N132").Value = -110110
$ws.Range("H134").Value = 11752.667
$ws.Range("I134").Value = 3448.111
$ws.Range("K134").Value = 10344.333
$ws.Range("M134").Value = -7809.332999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 45500
$ws.Range("J9").Value = 45500
$ws.Range("L9").Value = 45500
$ws.Range("N9").Value = -45836
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("H99").Value = 3783.2727
$ws.Range("J99").Value = 3076
$ws.Range("L99").Value = 3076
$ws.Range("N99").Value = -6072
$ws.Range("H126").Value = 3783.2727
$ws.Range("J126").Value = 3076
$ws.Range("L126").Value = 9228
$ws.Range("N126").Value = -14168
$ws.Range("M13").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 2600
$ws.Range("I63").Value = 2200
$ws.Range("K63").Value = 6600
$ws.Range("M63").Value = -5851
$ws.Range("H66").Value = 2600
$ws.Range("I66").Value = 2200
$ws.Range("K66").Value = 19800
$ws.Range("M66").Value = -16056
$ws.Range("H68").Value = 5087.222
$ws.Range("I68").Value = 288
$ws.Range("J68").Value = 5687.125
$ws.Range("K68").Value = 864
$ws.Range("L68").Value = 17061.375
$ws.Range("M68").Value = -53
$ws.Range("N68").Value = -18683.375
$ws.Range("H69").Value = 3000
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("H71").Value = 5087.222
$ws.Range("I71").Value = 288
$ws.Range("J71").Value = 5687.125
$ws.Range("K71").Value = 2592
$ws.Range("L71").Value = 51184.125
$ws.Range("M71").Value = 1464
$ws.Range("N71").Value = -59296.125
$ws.Range("H72").Value = 3000
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("H74").Value = 14699
$ws.Range("I74").Value = 14699
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 44097
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -43036
$ws.Range("H77").Value = 14699
$ws.Range("I77").Value = 14699
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 132291
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -126987
$ws.Range("H131").Value = 3092
$ws.Range("J131").Value = 3439.85
$ws.Range("L131").Value = 10319.55
$ws.Range("N131").Value = -20399.55
$ws.Range("M69").ClearContents()
$ws.Range("M72").ClearContents()
$ws.Range("N74").ClearContents()
$ws.Range("N77").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 23329
$ws.Range("J20").Value = 24993.5
$ws.Range("L20").Value = 24993.5
$ws.Range("N20").Value = -25483.5
$ws.Range("H24").Value = 157428
$ws.Range("J24").Value = 16999.334
$ws.Range("L24").Value = 16999.334
$ws.Range("N24").Value = -17345.334
$ws.Range("H80").Value = 2000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 2000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 2000
$ws.Range("N80").Value = -3996
$ws.Range("H83").Value = 2000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 2000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 10000
$ws.Range("N83").Value = -19984
$ws.Range("H113").Value = 3629.3845
$ws.Range("I113").Value = 2676.125
$ws.Range("K113").Value = 2676.125
$ws.Range("M113").Value = -506.125
$ws.Range("M80").ClearContents()
$ws.Range("M83").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1782.5667
$ws.Range("I55").Value = 2216
$ws.Range("K55").Value = 2216
$ws.Range("M55").Value = -2043
$ws.Range("H133").Value = 72500
$ws.Range("J133").Value = 72500
$ws.Range("L133").Value = 72500
$ws.Range("N133").Value = -77560

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("H52").Value = 12021
$ws.Range("I52").Value = 12021
$ws.Range("K52").Value = 12021
$ws.Range("M52").Value = -11795
$ws.Range("H81").Value = 2694.5557
$ws.Range("I81").Value = 2069.2307
$ws.Range("J81").Value = 4320.4
$ws.Range("K81").Value = 4138.4614
$ws.Range("L81").Value = 8640.799999999999
$ws.Range("M81").Value = -3077.4614
$ws.Range("N81").Value = -10762.8
$ws.Range("H84").Value = 2694.5557
$ws.Range("I84").Value = 2069.2307
$ws.Range("J84").Value = 4320.4
$ws.Range("K84").Value = 20692.307
$ws.Range("L84").Value = 43204
$ws.Range("M84").Value = -15388.307
$ws.Range("N84").Value = -53812
$ws.Range("H136").Value = 4279.9
$ws.Range("I136").Value = 3860
$ws.Range("J136").Value = 4699.8
$ws.Range("K136").Value = 11580
$ws.Range("L136").Value = 14099.4
$ws.Range("M136").Value = -9030
$ws.Range("N136").Value = -19199.4
$ws.Range("M31").ClearContents()
$ws.Range("N44").ClearContents()
$ws.Range("M51").ClearContents()
